# Target the "normal" sheet (the tab that was active / edited in the diff).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("normal")

# --- Probability-table tweaks (row 4 / row 5) ---
$ws.Range("G4").Value = 0.1
$ws.Range("I4").Value = 0.1

$ws.Range("F5").Value = 0.1
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 0.1
$ws.Range("I5").Value = 0

# --- View state: scroll so row 4 is the top visible row, select row 9 ---
$ws.Activate()
try { $excel.ActiveWindow.ScrollRow = 4 } catch {}
try { $excel.ActiveWindow.Panes.Item(1).TopLeftCell = $ws.Range("A4") } catch {}
$ws.Range("A9:XFD9").Select()
